# ---------------------------------------------------------------------------
# Edit: B1--and-B2-PowerPoint.pptx
#   1) Slide 5 table -> change table style (tableStyleId) to
#      {D60DA202-7103-4439-9964-AD1765470DDE}
#   2) The deck's theme colour scheme is swapped from the custom
#      "Red Violet"/"Integral" palette to the stock "Office" palette
#      (the commit swaps ppt/theme/theme1.xml <-> ppt/theme/theme2.xml;
#      since fontScheme/fmtScheme are identical between the two theme
#      parts, the only observable difference is the clrScheme, which we
#      reproduce here via the theme colour scheme object model).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- helper: "RRGGBB" hex string -> COM RGB long (0x00BBGGRR) -------------
function HexToRgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# ---------------------------------------------------------------------------
# 1) Table style on slide 5 (shape 2 is the 3-column table)
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{D60DA202-7103-4439-9964-AD1765470DDE}")
}

# ---------------------------------------------------------------------------
# 2) Theme colours -> stock "Office" palette
#    (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink in that order)
# ---------------------------------------------------------------------------
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = HexToRgbLong($officeColors[$i - 1])
}
